$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 176457951.9162543
    3 = 299025967.0485092
    4 = 387764540.3603849
    5 = 450452669.8438839
    6 = 493503532.9210405
    7 = 521775846.5592073
    8 = 539013044.5301551
    9 = 548092107.668923
    10 = 551218494.63842
    11 = 550076642.2618529
    12 = 545947163.2236661
    13 = 539799075.2234254
    14 = 532360431.3512471
    15 = 524173546.3835164
    16 = 515638975.8569094
    17 = 507047121.2556146
    18 = 498605368.38407
    19 = 490457923.8681079
    20 = 482701208.1152577
    21 = 475395156.8577317
    22 = 468573397.5775556
    23 = 462250059.8493025
    24 = 456425339.6149201
    25 = 451088589.7846206
    26 = 446222651.5220402
    27 = 441806171.1435871
    28 = 437814746.1119272
    29 = 434223109.8687769
    30 = 431004295.1508185
    31 = 428132409.2943134
    32 = 425582432.9840421
    33 = 423329516.2197678
    34 = 421351454.9057027
    35 = 419625581.6405653
    36 = 418132041.709063
    37 = 416851519.7363951
    38 = 415767662.3020656
    39 = 414863048.7489735
    40 = 414123986.309943
    41 = 413536515.9562929
    42 = 413088447.9028727
    43 = 412768106.2184888
    44 = 412566006.4589693
    45 = 412472108.5052615
    46 = 412477858.5543067
    47 = 412576282.6792486
    48 = 412759470.6436436
    49 = 413022121.3915059
    50 = 413358140.4965709
    51 = 413762482.4607639
    52 = 414229516.5931035
    53 = 414755968.4239311
    54 = 415338571.6120818
    55 = 415972866.3446266
    56 = 416656313.7065665
    57 = 417386167.9880645
    58 = 418159974.4036821
    59 = 418974754.8515918
    60 = 419829170.9755239
    61 = 420720852.4462495
    62 = 421648814.7359372
    63 = 422611708.927172
    64 = 423608223.9327038
    65 = 424636467.7553445
    66 = 425696278.871001
    67 = 426786074.136685
    68 = 427905958.6137171
    69 = 429054011.1814067
    70 = 430229911.7087038
    71 = 431433165.8240107
    72 = 432664340.1367074
    73 = 433921833.3599676
    74 = 435205640.3568087
    75 = 436515656.0983549
    76 = 437851493.7040353
    77 = 439212820.1531971
    78 = 440599422.3752252
    79 = 442011063.5524056
    80 = 443447991.5123991
    81 = 444909454.6380911
    82 = 446395544.7299273
    83 = 447906634.8050544
    84 = 449442090.6922197
    85 = 451001894.1783384
    86 = 452587010.978581
    87 = 454197334.7154202
    88 = 455831884.8775792
    89 = 457490854.9534469
    90 = 459174633.1963089
    91 = 460883060.0734919
    92 = 462616110.5870689
    93 = 464374571.0869839
    94 = 466157838.2171621
    95 = 467965785.8291762
    96 = 469797513.1694354
    97 = 471654710.3175163
    98 = 473536494.4916843
    99 = 475442985.1569877
    100 = 477374105.9818461
    101 = 479330091.3508222
    102 = 481310331.861344
    103 = 483315317.4137437
    104 = 485344596.9995273
    105 = 487397679.9636538
    106 = 489475409.6015839
    107 = 491576931.1183516
    108 = 493701771.6103851
    109 = 495849257.1985958
    110 = 498019733.3222932
    111 = 500213906.9507883
    112 = 502430983.3973305
    113 = 504669070.1116637
    114 = 506929664.4962125
    115 = 509212060.9986465
    116 = 511515780.2451865
    117 = 513840737.2478815
    118 = 516185984.8176445
    119 = 518552154.4726911
    120 = 520938221.8249016
    121 = 523343420.5276664
    122 = 525766984.7309943
    123 = 528209449.5818053
    124 = 530669261.1211135
    125 = 533146854.0029674
    126 = 535640194.0273566
    127 = 538149448.5615404
    128 = 540674451.3420504
    129 = 543214624.9452326
    130 = 545769154.5474751
    131 = 548337072.0779687
    132 = 550917973.5203508
    133 = 553512649.3146377
    134 = 556119535.0484209
    135 = 558737902.6524844
    136 = 561367362.3105721
    137 = 564007588.3944304
    138 = 566657332.6379468
    139 = 569316913.5714279
    140 = 571985054.2705294
    141 = 574661874.0533047
    142 = 577347537.2589236
    143 = 580041779.4814882
    144 = 582743630.0038074
    145 = 585453810.3281393
    146 = 588171576.6691861
    147 = 590897187.1987077
    148 = 593631035.273096
    149 = 596371703.1228534
    150 = 599120504.6910335
    151 = 601877586.7182524
    152 = 604643274.495908
    153 = 607418706.8039639
    154 = 610203718.7511187
    155 = 612999358.6074624
    156 = 615806070.8515373
    157 = 618624521.2447153
    158 = 621456146.266559
    159 = 624301629.9625001
    160 = 627162425.890736
    161 = 630039882.3513825
    162 = 632934356.9967171
    163 = 635849098.7305893
    164 = 638784862.2847214
    165 = 641742482.0348504
    166 = 644725819.9458139
    167 = 647734287.9686497
    168 = 650771703.7271391
    169 = 653839022.9449497
    170 = 656939669.4470539
    171 = 660074795.3635426
    172 = 663246547.478405
    173 = 666457320.4314189
    174 = 669709449.1104873
    175 = 673006822.6657867
    176 = 676350739.7155666
    177 = 679743909.7535046
    178 = 683190219.321215
    179 = 686690663.4082291
    180 = 690249571.4971693
    181 = 693869308.0067152
    182 = 697552350.2227504
    183 = 701301803.6859961
    184 = 705121258.8353307
    185 = 709014046.3522103
    186 = 712982558.1323545
    187 = 717030177.1548752
    188 = 721160266.6619897
    189 = 725375746.0260571
    190 = 729678434.64591
    191 = 734075421.2923222
    192 = 738568383.6500947
    193 = 743160967.6170096
    194 = 747855622.3948622
    195 = 752657535.7208365
    196 = 757570006.042662
    197 = 762597151.6234024
    198 = 767742026.9118183
    199 = 773009365.591303
    200 = 778402664.7883159
    201 = 783926911.249671
    202 = 789585589.2733015
    203 = 795382371.6121786
    204 = 801323482.4720855
    205 = 807412194.296483
    206 = 813653659.6441444
    207 = 820053354.6124957
    208 = 826614213.5386038
    209 = 833342489.9879069
    210 = 840242819.8908347
    211 = 847320704.4161122
    212 = 854581575.1523492
    213 = 862029945.1327573
    214 = 869672536.6917185
    215 = 877514637.5087852
    216 = 885562321.5754826
    217 = 893821508.8273773
    218 = 902298563.9078224
    219 = 910998452.4936625
    220 = 919930152.7007653
    221 = 929098140.9319403
    222 = 938509661.9031746
    223 = 948171557.8569007
    224 = 958090796.4427094
    225 = 968274225.3173132
    226 = 978730502.2457995
    227 = 989465871.8227062
    228 = 1000487852.780719
    229 = 1011805970.587862
    230 = 1023426075.220151
    231 = 1035357698.876846
    232 = 1047608755.683886
    233 = 1060187321.811259
    234 = 1073103452.944369
    235 = 1086364928.015782
    236 = 1099980922.256836
    237 = 1113960420.460771
    238 = 1128313398.937002
    239 = 1143049667.916685
    240 = 1158177579.583402
    241 = 1173708022.41945
    242 = 1189652144.905309
    243 = 1206019544.495125
    244 = 1222819684.257553
    245 = 1240063856.72824
    246 = 1257763562.373686
    247 = 1275931508.633592
    248 = 1294577783.051069
    249 = 1313712139.349321
    250 = 1333349660.391406
    251 = 1353499295.12593
    252 = 1374178237.39651
    253 = 1395394362.053137
    254 = 1417162920.509463
    255 = 1439495226.142487
    256 = 1462406707.834435
    257 = 1485908034.920981
    258 = 1510015482.280327
    259 = 1534741671.632022
    260 = 1560102279.012724
    261 = 1586111362.865478
    262 = 1612784092.23097
    263 = 1640136354.272816
    264 = 1668182906.56835
    265 = 1696939612.13547
    266 = 1726422749.812349
    267 = 1756649365.036398
    268 = 1787635806.158347
    269 = 1819399930.776184
    270 = 1851958321.761299
    271 = 1885329132.137106
    272 = 1919532337.069093
    273 = 1954585337.662289
    274 = 1990507097.540207
    275 = 2027318166.897348
    276 = 2065037461.240351
    277 = 2103685628.385842
    278 = 2143282305.135448
    279 = 2183851072.946911
    280 = 2225413344.823607
    281 = 2267990794.331093
    282 = 2311605942.896164
    283 = 2356282983.223853
    284 = 2402043918.052493
    285 = 2448915508.89982
    286 = 2496923500.478505
    287 = 2546089488.51217
    288 = 2596444090.239007
    289 = 2648011289.784314
    290 = 2700820446.309979
    291 = 2754897491.071284
    292 = 2810275132.994972
    293 = 2866979238.595462
    294 = 2925042736.993287
    295 = 2984496882.735734
    296 = 3045370666.285125
    297 = 3107700655.394271
    298 = 3171518620.472177
    299 = 3236859399.959856
    300 = 3303760612.959435
    301 = 3372256028.132004
    302 = 3442381941.584273
    303 = 3512910734.50428
    304 = 3584555905.983273
    305 = 3657344107.275531
    306 = 3731305075.702945
    307 = 3806468679.800125
    308 = 3882862312.164623
    309 = 3960516149.154913
    310 = 4039464208.048369
    311 = 4119735917.129042
    312 = 4201366964.401301
    313 = 4284387866.502783
    314 = 4368835935.102097
    315 = 4454747944.437763
    316 = 4542158202.830487
    317 = 4631103275.857863
    318 = 4721626165.148641
    319 = 4813765722.541456
    320 = 4907562209.034718
    321 = 5003059906.442801
    322 = 5100300517.066927
    323 = 5199329446.016441
    324 = 5300195727.460501
    325 = 5402948152.083224
    326 = 5507634380.946594
    327 = 5614302493.463202
    328 = 5723011387.660942
    329 = 5833812941.289371
    330 = 5946761784.844069
    331 = 6061918208.294804
    332 = 6179338008.449124
    333 = 6299086486.289141
    334 = 6421227366.715107
    335 = 6545824971.192794
    336 = 6672944804.970401
    337 = 6802660583.255257
    338 = 6935044074.298243
    339 = 7070169718.629631
    340 = 7208112923.549122
    341 = 7348957614.593982
    342 = 7492785161.435617
    343 = 7639678885.241313
    344 = 7789726954.888515
    345 = 7943022958.694237
    346 = 8099656347.415995
    347 = 8259730043.379133
    348 = 8423341106.135232
    349 = 8590594357.845703
    350 = 8761598166.30127
    351 = 8936465953.29496
    352 = 9115311567.803825
    353 = 9298257878.880634
    354 = 9485424377.676859
    355 = 9676942909.445187
    356 = 9872947729.406685
    357 = 10073574224.81887
    358 = 10279325116.51797
    359 = 10490162218.28854
    360 = 10706244114.86433
    361 = 10927735692.70711
    362 = 11154809250.26011
    363 = 11387644623.97968
    364 = 11626431517.58892
    365 = 11871362080.94988
    366 = 12122636003.54384
    367 = 12380466953.59603
    368 = 12645074370.33446
    369 = 12916687592.5314
    370 = 13195546158.65542
    371 = 13481900415.61717
    372 = 13776011578.64161
    373 = 14075158479.95717
    374 = 14381201267.90539
    375 = 14694394198.27159
    376 = 15015225243.30906
    377 = 15343762140.17858
    378 = 15680300069.25786
    379 = 16025150077.5308
    380 = 16378643589.37771
    381 = 16741129531.55394
    382 = 17112978453.47021
    383 = 17500697141.71965
    384 = 17901530639.41035
    385 = 18316042930.92529
    386 = 18744835079.62518
    387 = 19188556256.72184
    388 = 19679847365.55902
    389 = 20202623789.5269
    390 = 20753856019.99387
    391 = 21321990141.67669
    392 = 21889416124.71362
    393 = 22403262617.79168
    394 = 22865937546.8485
    395 = 23268908587.33424
    396 = 23610238225.26783
    397 = 23893262830.29799
    398 = 24124926430.82415
    399 = 24314241765.25021
    400 = 24471074485.25926
    401 = 24605301034.69932
    402 = 24726288567.50959
    403 = 24741768276.72356
    404 = 24707014587.29932
    405 = 24633961654.30794
    406 = 24533230784.36732
    407 = 24414124793.88937
    408 = 24284733419.43191
    409 = 24152177870.23462
    410 = 24022941412.5499
    411 = 23903461050.73541
    412 = 23801118914.53946
    413 = 23752788516.19056
    414 = 23771402499.44581
    415 = 23819906859.53962
    416 = 23891542949.0889
    417 = 23980750688.00551
    418 = 24082984290.53928
    419 = 24192350148.10149
    420 = 24281844331.04657
    421 = 24287632138.12484
    422 = 24136129767.32555
    423 = 23779004623.43476
    424 = 23206165298.69995
    425 = 22439902515.70296
    426 = 21521891600.47131
    427 = 20500933409.36682
    428 = 19424321223.95302
    429 = 18332981692.67078
    430 = 17259495439.60363
    431 = 16227936206.26147
    432 = 15254731803.80456
    433 = 14349927377.02803
    434 = 13518562188.15435
    435 = 12761907526.90144
    436 = 12078549434.05006
    437 = 11465268395.7449
    438 = 10917706829.31154
    439 = 10430877159.91087
    440 = 9999531429.530464
    441 = 9618415687.526512
    442 = 9282455983.073681
    443 = 8986844186.591637
    444 = 8727123561.347092
    445 = 8499203168.85235
    446 = 8299380849.412382
    447 = 8124319685.007588
    448 = 7971042846.380408
    449 = 7836885206.220637
    450 = 7719501188.616909
    451 = 7616822589.862211
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
